$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "q" to every QID2 value in column A (rows 2..73), leaving the
# header row (row 1) untouched.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 73) { $lastRow = 73 }
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2) {
        $cell.Value2 = "q" + $cell.Value2
    }
}

# Update the view: clear the frozen/scrolled "topLeftCell" and move the
# selection from F1:F1048576 to F15.
$ws.Range("F15").Select()
